$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two JIRA URL rows that survive to the new content.
$ws.Range("A2").Value = "https://jira.jnj.com/browse/JCVZ-998"
$ws.Range("A3").Value = "https://jira.jnj.com/browse/JCVZ-1035"

# Drop the now-obsolete rows (old rows 4-11).
$ws.Range("A4:A11").EntireRow.Delete()

# Clear every hyperlink on the sheet (scoped Range.Hyperlinks.Delete affects
# the whole sheet in this runtime), then re-add only the one that should
# remain: A3 pointing at the new JCVZ-1035 ticket.
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A3"), "https://jira.jnj.com/browse/JCVZ-1035")

# Re-apply the Hyperlink cell style so A3 keeps using the same shared style
# as A2 instead of a freshly minted one.
$ws.Range("A3").Style = "Hyperlink"

# Move the active selection to match the saved view state.
[void]$ws.Range("E6").Select()
